# Apply scheduled market-price / profit recalculation updates to each Leve profit sheet.
# Values correspond to cells in columns H:N (currentAveragePrice*, LevePrice*, LeveProfit*)
# for the specific rows that were refreshed by the pricing runner.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(19, 8).Value = 487.875  # H19
$ws.Cells.Item(19, 9).Value = 475.25  # I19
$ws.Cells.Item(19, 10).Value = 500.5  # J19
$ws.Cells.Item(19, 11).Value = 475.25  # K19
$ws.Cells.Item(19, 12).Value = 500.5  # L19
$ws.Cells.Item(19, 13).Value = -300.25  # M19
$ws.Cells.Item(19, 14).Value = -850.5  # N19
$ws.Cells.Item(112, 8).Value = 6612.6553  # H112
$ws.Cells.Item(112, 10).Value = 7856.9585  # J112
$ws.Cells.Item(112, 12).Value = 23570.8755  # L112
$ws.Cells.Item(112, 14).Value = -25786.8755  # N112
$ws.Cells.Item(121, 8).Value = 579.8387  # H121
$ws.Cells.Item(121, 10).Value = 539.1667  # J121
$ws.Cells.Item(121, 12).Value = 1617.5001  # L121
$ws.Cells.Item(121, 14).Value = -5111.5001  # N121

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(37, 8).Value = 6789.4  # H37
$ws.Cells.Item(37, 9).Value = 1034  # I37
$ws.Cells.Item(37, 10).Value = 8228.25  # J37
$ws.Cells.Item(37, 11).Value = 1034  # K37
$ws.Cells.Item(37, 12).Value = 8228.25  # L37
$ws.Cells.Item(37, 13).Value = -761  # M37
$ws.Cells.Item(37, 14).Value = -8774.25  # N37
$ws.Cells.Item(45, 8).Value = 2243.9773  # H45
$ws.Cells.Item(45, 9).Value = 2204.7334  # I45
$ws.Cells.Item(45, 10).Value = 2328.0715  # J45
$ws.Cells.Item(45, 11).Value = 2204.7334  # K45
$ws.Cells.Item(45, 12).Value = 2328.0715  # L45
$ws.Cells.Item(45, 13).Value = -1827.7334  # M45
$ws.Cells.Item(45, 14).Value = -3082.0715  # N45
$ws.Cells.Item(61, 8).Value = 4545.8184  # H61
$ws.Cells.Item(61, 10).Value = 2704.6667  # J61
$ws.Cells.Item(61, 12).Value = 2704.6667  # L61
$ws.Cells.Item(61, 14).Value = -3128.6667  # N61
$ws.Cells.Item(97, 8).Value = 3149.8333  # H97
$ws.Cells.Item(97, 9).Value = 2705.0952  # I97
$ws.Cells.Item(97, 10).Value = 6263  # J97
$ws.Cells.Item(97, 11).Value = 2705.0952  # K97
$ws.Cells.Item(97, 12).Value = 6263  # L97
$ws.Cells.Item(97, 13).Value = -2209.0952  # M97
$ws.Cells.Item(97, 14).Value = -7255  # N97
$ws.Cells.Item(132, 8).Value = 3070.5356  # H132
$ws.Cells.Item(132, 9).Value = 2716.611  # I132
$ws.Cells.Item(132, 10).Value = 3707.6  # J132
$ws.Cells.Item(132, 11).Value = 8149.833  # K132
$ws.Cells.Item(132, 12).Value = 11122.8  # L132
$ws.Cells.Item(132, 13).Value = -5619.833  # M132
$ws.Cells.Item(132, 14).Value = -16182.8  # N132
$ws.Cells.Item(136, 8).Value = 4545.8184  # H136
$ws.Cells.Item(136, 10).Value = 2704.6667  # J136
$ws.Cells.Item(136, 12).Value = 8114.000100000001  # L136
$ws.Cells.Item(136, 14).Value = -13214.0001  # N136

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(57, 8).Value = 39990  # H57
$ws.Cells.Item(57, 10).Value = 39990  # J57
$ws.Cells.Item(57, 12).Value = 39990  # L57
$ws.Cells.Item(57, 14).Value = -41430  # N57
$ws.Cells.Item(134, 8).Value = 5039.5312  # H134
$ws.Cells.Item(134, 9).Value = 497.7857  # I134
$ws.Cells.Item(134, 10).Value = 36831.75  # J134
$ws.Cells.Item(134, 11).Value = 1493.3571  # K134
$ws.Cells.Item(134, 12).Value = 110495.25  # L134
$ws.Cells.Item(134, 13).Value = 1041.6429  # M134
$ws.Cells.Item(134, 14).Value = -115565.25  # N134
$ws.Cells.Item(136, 8).Value = 39990  # H136
$ws.Cells.Item(136, 10).Value = 39990  # J136
$ws.Cells.Item(136, 12).Value = 39990  # L136
$ws.Cells.Item(136, 14).Value = -50190  # N136
$ws.Cells.Item(138, 8).Value = 50750  # H138
$ws.Cells.Item(138, 10).Value = 50750  # J138
$ws.Cells.Item(138, 12).Value = 50750  # L138
$ws.Cells.Item(138, 14).Value = -61030  # N138
$ws.Cells.Item(140, 8).Value = 86690  # H140
$ws.Cells.Item(140, 10).Value = 86690  # J140
$ws.Cells.Item(140, 12).Value = 86690  # L140
$ws.Cells.Item(140, 14).Value = -97050  # N140

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(51, 8).Value = 9498.5  # H51
$ws.Cells.Item(51, 10).Value = 9498.5  # J51
$ws.Cells.Item(51, 12).Value = 9498.5  # L51
$ws.Cells.Item(51, 14).Value = -10970.5  # N51
$ws.Cells.Item(58, 8).Value = 1550.55  # H58
$ws.Cells.Item(58, 9).Value = 924.8889  # I58
$ws.Cells.Item(58, 10).Value = 2062.4546  # J58
$ws.Cells.Item(58, 11).Value = 924.8889  # K58
$ws.Cells.Item(58, 12).Value = 2062.4546  # L58
$ws.Cells.Item(58, 13).Value = -721.8889  # M58
$ws.Cells.Item(58, 14).Value = -2468.4546  # N58
$ws.Cells.Item(61, 8).Value = 9498.5  # H61
$ws.Cells.Item(61, 10).Value = 9498.5  # J61
$ws.Cells.Item(61, 12).Value = 9498.5  # L61
$ws.Cells.Item(61, 14).Value = -10194.5  # N61
$ws.Cells.Item(74, 8).Value = 17437.715  # H74
$ws.Cells.Item(74, 10).Value = 17437.715  # J74
$ws.Cells.Item(74, 12).Value = 17437.715  # L74
$ws.Cells.Item(74, 14).Value = -19185.715  # N74
$ws.Cells.Item(77, 8).Value = 17437.715  # H77
$ws.Cells.Item(77, 10).Value = 17437.715  # J77
$ws.Cells.Item(77, 12).Value = 52313.145  # L77
$ws.Cells.Item(77, 14).Value = -61049.145  # N77
$ws.Cells.Item(136, 8).Value = 1550.55  # H136
$ws.Cells.Item(136, 9).Value = 924.8889  # I136
$ws.Cells.Item(136, 10).Value = 2062.4546  # J136
$ws.Cells.Item(136, 11).Value = 2774.6667  # K136
$ws.Cells.Item(136, 12).Value = 6187.3638  # L136
$ws.Cells.Item(136, 13).Value = -224.6667000000002  # M136
$ws.Cells.Item(136, 14).Value = -11287.3638  # N136

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 8).Value = 188.5  # H2
$ws.Cells.Item(2, 9).Value = 40.214287  # I2
$ws.Cells.Item(2, 10).Value = 336.7857  # J2
$ws.Cells.Item(2, 11).Value = 241.285722  # K2
$ws.Cells.Item(2, 12).Value = 2020.7142  # L2
$ws.Cells.Item(2, 13).Value = -128.285722  # M2
$ws.Cells.Item(2, 14).Value = -2246.7142  # N2
$ws.Cells.Item(9, 8).Value = 25103350  # H9
$ws.Cells.Item(9, 9).Value = 266667  # I9
$ws.Cells.Item(9, 10).Value = 40005360  # J9
$ws.Cells.Item(9, 11).Value = 800001  # K9
$ws.Cells.Item(9, 12).Value = 120016080  # L9
$ws.Cells.Item(9, 13).Value = -799777  # M9
$ws.Cells.Item(9, 14).Value = -120016528  # N9
$ws.Cells.Item(12, 8).Value = 381.5  # H12
$ws.Cells.Item(12, 9).Value = 0  # I12
$ws.Cells.Item(12, 10).Value = 381.5  # J12
$ws.Cells.Item(12, 11).Value = 0  # K12
$ws.Cells.Item(12, 12).Value = 1144.5  # L12
$ws.Cells.Item(12, 13).ClearContents()  # M12
$ws.Cells.Item(12, 14).Value = -1490.5  # N12
$ws.Cells.Item(13, 8).Value = 206.66667  # H13
$ws.Cells.Item(13, 9).Value = 60  # I13
$ws.Cells.Item(13, 11).Value = 180  # K13
$ws.Cells.Item(13, 13).Value = -12  # M13
$ws.Cells.Item(16, 8).Value = 4153.6665  # H16
$ws.Cells.Item(16, 9).Value = 501  # I16
$ws.Cells.Item(16, 10).Value = 5980  # J16
$ws.Cells.Item(16, 11).Value = 1503  # K16
$ws.Cells.Item(16, 12).Value = 17940  # L16
$ws.Cells.Item(16, 13).Value = -1330  # M16
$ws.Cells.Item(16, 14).Value = -18286  # N16
$ws.Cells.Item(17, 8).Value = 283.75  # H17
$ws.Cells.Item(17, 9).Value = 94  # I17
$ws.Cells.Item(17, 10).Value = 600  # J17
$ws.Cells.Item(17, 11).Value = 282  # K17
$ws.Cells.Item(17, 12).Value = 1800  # L17
$ws.Cells.Item(17, 13).Value = -113  # M17
$ws.Cells.Item(17, 14).Value = -2138  # N17
$ws.Cells.Item(19, 8).Value = 500  # H19
$ws.Cells.Item(19, 9).Value = 500  # I19
$ws.Cells.Item(19, 11).Value = 1500  # K19
$ws.Cells.Item(19, 13).Value = -1326  # M19
$ws.Cells.Item(20, 8).Value = 1975  # H20
$ws.Cells.Item(20, 10).Value = 1975  # J20
$ws.Cells.Item(20, 12).Value = 5925  # L20
$ws.Cells.Item(20, 14).Value = -6379  # N20
$ws.Cells.Item(23, 8).Value = 472.05554  # H23
$ws.Cells.Item(23, 9).Value = 251.25  # I23
$ws.Cells.Item(23, 10).Value = 913.6667  # J23
$ws.Cells.Item(23, 11).Value = 753.75  # K23
$ws.Cells.Item(23, 12).Value = 2741.0001  # L23
$ws.Cells.Item(23, 13).Value = -518.75  # M23
$ws.Cells.Item(23, 14).Value = -3211.0001  # N23
$ws.Cells.Item(56, 8).Value = 3606775  # H56
$ws.Cells.Item(56, 9).Value = 3606775  # I56
$ws.Cells.Item(56, 11).Value = 3606775  # K56
$ws.Cells.Item(56, 13).Value = -3606245  # M56
$ws.Cells.Item(122, 8).Value = 501.4634  # H122
$ws.Cells.Item(122, 9).Value = 354.58823  # I122
$ws.Cells.Item(122, 10).Value = 1214.8572  # J122
$ws.Cells.Item(122, 11).Value = 3191.29407  # K122
$ws.Cells.Item(122, 12).Value = 10933.7148  # L122
$ws.Cells.Item(122, 13).Value = -741.2940699999999  # M122
$ws.Cells.Item(122, 14).Value = -15833.7148  # N122

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(113, 8).Value = 6472407.5  # H113
$ws.Cells.Item(113, 10).Value = 911051.0600000001  # J113
$ws.Cells.Item(113, 12).Value = 911051.0600000001  # L113
$ws.Cells.Item(113, 14).Value = -915391.0600000001  # N113
$ws.Cells.Item(132, 8).Value = 3491.6667  # H132
$ws.Cells.Item(132, 9).Value = 2923.111  # I132
$ws.Cells.Item(132, 10).Value = 4344.5  # J132
$ws.Cells.Item(132, 11).Value = 8769.332999999999  # K132
$ws.Cells.Item(132, 12).Value = 13033.5  # L132
$ws.Cells.Item(132, 13).Value = -6239.332999999999  # M132
$ws.Cells.Item(132, 14).Value = -18093.5  # N132
$ws.Cells.Item(138, 8).Value = 69866.664  # H138
$ws.Cells.Item(138, 10).Value = 69866.664  # J138
$ws.Cells.Item(138, 12).Value = 69866.664  # L138
$ws.Cells.Item(138, 14).Value = -80146.664  # N138
$ws.Cells.Item(140, 8).Value = 75776.336  # H140
$ws.Cells.Item(140, 10).Value = 75776.336  # J140
$ws.Cells.Item(140, 12).Value = 75776.336  # L140
$ws.Cells.Item(140, 14).Value = -86136.336  # N140
$ws.Cells.Item(141, 8).Value = 69900  # H141
$ws.Cells.Item(141, 10).Value = 69900  # J141
$ws.Cells.Item(141, 12).Value = 69900  # L141
$ws.Cells.Item(141, 14).Value = -80260  # N141

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(138, 8).Value = 60195.7  # H138
$ws.Cells.Item(138, 10).Value = 60195.7  # J138
$ws.Cells.Item(138, 12).Value = 60195.7  # L138
$ws.Cells.Item(138, 14).Value = -70475.7  # N138
$ws.Cells.Item(139, 8).Value = 79550  # H139
$ws.Cells.Item(139, 10).Value = 79550  # J139
$ws.Cells.Item(139, 12).Value = 79550  # L139
$ws.Cells.Item(139, 14).Value = -89830  # N139

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(126, 8).Value = 1618.85  # H126
$ws.Cells.Item(126, 9).Value = 892.8182  # I126
$ws.Cells.Item(126, 10).Value = 2506.2222  # J126
$ws.Cells.Item(126, 11).Value = 2678.4546  # K126
$ws.Cells.Item(126, 12).Value = 7518.6666  # L126
$ws.Cells.Item(126, 13).Value = -208.4546  # M126
$ws.Cells.Item(126, 14).Value = -12458.6666  # N126
$ws.Cells.Item(136, 8).Value = 1233.4482  # H136
$ws.Cells.Item(136, 9).Value = 716.5909  # I136
$ws.Cells.Item(136, 10).Value = 2857.8572  # J136
$ws.Cells.Item(136, 11).Value = 2149.7727  # K136
$ws.Cells.Item(136, 12).Value = 8573.571599999999  # L136
$ws.Cells.Item(136, 13).Value = 400.2273  # M136
$ws.Cells.Item(136, 14).Value = -13673.5716  # N136

